# "act tablas web jul25" — refresh the indicator table with Jul-2025 data
# (adds 2023 and 2022 data points, updates the historical series, and
# records the metadata "actualizacion" / last-update note).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Data sheet: insert the two newest years at the top and push the
# existing series down, then refresh every value with the latest
# published figures.
# ---------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Rows("2:3").Insert()

$wsData.Range("A2").NumberFormat = "@"
$wsData.Range("A2").Value = "2023"
$wsData.Range("A2").Style = "Normal"

$wsData.Range("A3").NumberFormat = "@"
$wsData.Range("A3").Value = "2022"
$wsData.Range("A3").Style = "Normal"

$wsData.Range("B2").Value = 726.3
$wsData.Range("B3").Value = 664.7
$wsData.Range("B4").Value = 628.8
$wsData.Range("B5").Value = 653.9
$wsData.Range("B6").Value = 670.5
$wsData.Range("B7").Value = 651.6
$wsData.Range("B8").Value = 652
$wsData.Range("B9").Value = 612.8
$wsData.Range("B10").Value = 531.7
$wsData.Range("B11").Value = 516.5
$wsData.Range("B12").Value = 478.2
$wsData.Range("B13").Value = 426.1
$wsData.Range("B14").Value = 386
$wsData.Range("B15").Value = 328.1
$wsData.Range("B16").Value = 297
$wsData.Range("B17").Value = 265.5
$wsData.Range("B18").Value = 221.4
$wsData.Range("B19").Value = 205
$wsData.Range("B20").Value = 176.1
$wsData.Range("B21").Value = 166.4
$wsData.Range("B22").Value = 174
$wsData.Range("B23").Value = 172.9
$wsData.Range("B24").Value = 191.3
$wsData.Range("B25").Value = 208.8
$wsData.Range("B26").Value = 202.6
$wsData.Range("B27").Value = 178.4
$wsData.Range("B28").Value = 168.3
$wsData.Range("B29").Value = 164.4
$wsData.Range("B30").Value = 155.8
$wsData.Range("B31").Value = 172.3
$wsData.Range("B32").Value = 122.2
$wsData.Range("B33").Value = 121.8
$wsData.Range("B34").Value = 105.3
$wsData.Range("B35").Value = 98.7

# ---------------------------------------------------------------
# Metadata sheet: blank header cell becomes a single space (matching
# the sibling cell), and a new "actualizacion" row is recorded ahead
# of the "cita" row.
# ---------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("A1").Value = " "

$wsMeta.Rows("9:9").Insert()
$wsMeta.Range("A9").Value = "actualizacion"
$wsMeta.Range("B9").Value = "Julio 2025"
